$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-06 15:18:07'
$ws.Range('H2').Value = '''89%'
$ws.Range('K2').Value = '8.2 MJ/m2'
$ws.Range('O2').Value = '-0.3 °C'
$ws.Range('E3').Value = '2026-02-06 15:18:09'
$ws.Range('H3').Value = '''71%'
$ws.Range('K3').Value = '11.5 MJ/m2'
$ws.Range('E4').Value = '2026-02-06 15:18:12'
$ws.Range('J4').Value = '996.6 hPa'
$ws.Range('K4').Value = '10.9 MJ/m2'
$ws.Range('O4').Value = '13.4 °C'
$ws.Range('E5').Value = '2026-02-06 15:18:14'
$ws.Range('J5').Value = '996.9 hPa'
$ws.Range('K5').Value = '10.0 MJ/m2'
$ws.Range('O5').Value = '10.7 °C'
$ws.Range('E6').Value = '2026-02-06 15:18:17'
$ws.Range('J6').Value = '998.1 hPa'
$ws.Range('K6').Value = '9.1 MJ/m2'
$ws.Range('E7').Value = '2026-02-06 15:18:19'
$ws.Range('J7').Value = '997.7 hPa'
$ws.Range('K7').Value = '11.4 MJ/m2'
$ws.Range('O7').Value = '11.5 °C'
$ws.Range('E8').Value = '2026-02-06 15:18:22'
$ws.Range('H8').Value = '''77%'
$ws.Range('K8').Value = '11.2 MJ/m2'
$ws.Range('O8').Value = '9.8 °C'
$ws.Range('E9').Value = '2026-02-06 15:18:24'
$ws.Range('H9').Value = '''86%'
$ws.Range('O9').Value = '4.4 °C'
$ws.Range('E10').Value = '2026-02-06 15:18:27'
$ws.Range('H10').Value = '''88%'
$ws.Range('O10').Value = '8.6 °C'
$ws.Range('E11').Value = '2026-02-06 15:18:29'
$ws.Range('H11').Value = '''78%'
$ws.Range('K11').Value = '8.2 MJ/m2'
$ws.Range('O11').Value = '5.0 °C'
$ws.Range('E12').Value = '2026-02-06 15:18:32'
$ws.Range('K12').Value = '11.3 MJ/m2'
$ws.Range('O12').Value = '14.3 °C'
$ws.Range('E13').Value = '2026-02-06 15:18:35'
$ws.Range('H13').Value = '''77%'
$ws.Range('O13').Value = '9.8 °C'
$ws.Range('E14').Value = '2026-02-06 15:18:36'
$ws.Range('K14').Value = '7.1 MJ/m2'
$ws.Range('O14').Value = '-4.0 °C'
$ws.Range('E15').Value = '2026-02-06 15:18:39'
$ws.Range('H15').Value = '''73%'
$ws.Range('J15').Value = '997.0 hPa'
$ws.Range('K15').Value = '11.1 MJ/m2'
$ws.Range('O15').Value = '10.0 °C'
$ws.Range('E16').Value = '2026-02-06 15:18:42'
$ws.Range('K16').Value = '8.9 MJ/m2'
$ws.Range('M16').Value = '11.4 °C 14:59 TU'
$ws.Range('O16').Value = '5.7 °C'
$ws.Range('E17').Value = '2026-02-06 15:18:45'
$ws.Range('H17').Value = '''86%'
$ws.Range('K17').Value = '9.8 MJ/m2'
$ws.Range('L17').Value = '26.3 km/h - 236º 14:41 TU'
$ws.Range('O17').Value = '5.5 °C'
$ws.Range('E18').Value = '2026-02-06 15:18:47'
$ws.Range('K18').Value = '5.5 MJ/m2'
$ws.Range('O18').Value = '-4.3 °C'
$ws.Range('E19').Value = '2026-02-06 15:18:50'
$ws.Range('H19').Value = '''78%'
$ws.Range('K19').Value = '11.0 MJ/m2'
$ws.Range('O19').Value = '9.5 °C'
$ws.Range('E20').Value = '2026-02-06 15:18:53'
$ws.Range('K20').Value = '11.1 MJ/m2'
$ws.Range('E21').Value = '2026-02-06 15:18:56'
$ws.Range('H21').Value = '''75%'
$ws.Range('J21').Value = '997.3 hPa'
$ws.Range('K21').Value = '9.7 MJ/m2'
$ws.Range('O21').Value = '8.1 °C'
$ws.Range('E22').Value = '2026-02-06 15:18:59'
$ws.Range('H22').Value = '''76%'
$ws.Range('K22').Value = '10.8 MJ/m2'
$ws.Range('O22').Value = '10.4 °C'
$ws.Range('E23').Value = '2026-02-06 15:19:01'
$ws.Range('J23').Value = '997.0 hPa'
$ws.Range('K23').Value = '8.6 MJ/m2'
$ws.Range('O23').Value = '9.8 °C'
$ws.Range('E24').Value = '2026-02-06 15:19:04'
$ws.Range('H24').Value = '''63%'
$ws.Range('J24').Value = '996.4 hPa'
$ws.Range('K24').Value = '11.2 MJ/m2'
$ws.Range('O24').Value = '13.1 °C'
$ws.Range('E25').Value = '2026-02-06 15:19:07'
$ws.Range('H25').Value = '''81%'
$ws.Range('K25').Value = '8.7 MJ/m2'
$ws.Range('O25').Value = '4.0 °C'
$ws.Range('E26').Value = '2026-02-06 15:19:09'
$ws.Range('H26').Value = '''78%'
$ws.Range('K26').Value = '7.6 MJ/m2'
$ws.Range('O26').Value = '-0.9 °C'
$ws.Range('E27').Value = '2026-02-06 15:19:11'
$ws.Range('H27').Value = '''83%'
$ws.Range('J27').Value = '997.0 hPa'
$ws.Range('K27').Value = '9.7 MJ/m2'
$ws.Range('O27').Value = '10.4 °C'
$ws.Range('E28').Value = '2026-02-06 15:19:14'
$ws.Range('J28').Value = '999.4 hPa'
$ws.Range('O28').Value = '4.3 °C'
$ws.Range('E29').Value = '2026-02-06 15:19:17'
$ws.Range('K29').Value = '11.4 MJ/m2'
$ws.Range('O29').Value = '12.5 °C'
$ws.Range('E30').Value = '2026-02-06 15:19:20'
$ws.Range('H30').Value = '''75%'
$ws.Range('K30').Value = '8.6 MJ/m2'
$ws.Range('L30').Value = '43.6 km/h - 223º 14:50 TU'
$ws.Range('E31').Value = '2026-02-06 15:19:22'
$ws.Range('H31').Value = '''86%'
$ws.Range('O31').Value = '6.9 °C'
$ws.Range('E32').Value = '2026-02-06 15:19:25'
$ws.Range('H32').Value = '''47%'
$ws.Range('K32').Value = '11.4 MJ/m2'
$ws.Range('E33').Value = '2026-02-06 15:19:28'
$ws.Range('H33').Value = '''84%'
$ws.Range('O33').Value = '9.8 °C'
$ws.Range('E34').Value = '2026-02-06 15:19:30'
$ws.Range('K34').Value = '11.1 MJ/m2'
$ws.Range('O34').Value = '8.4 °C'
$ws.Range('E35').Value = '2026-02-06 15:19:33'
$ws.Range('K35').Value = '8.7 MJ/m2'
$ws.Range('O35').Value = '-2.1 °C'
$ws.Range('E36').Value = '2026-02-06 15:19:35'
$ws.Range('H36').Value = '''60%'
$ws.Range('K36').Value = '10.9 MJ/m2'
$ws.Range('O36').Value = '13.2 °C'
